# Auto-generated edit script: updates market-price derived columns (H-N)
# across 37 rows spanning all 8 class sheets, per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 777.55554
$ws.Cells.Item(43, 9).Value = 622.3333
$ws.Cells.Item(43, 10).Value = 855.1667
$ws.Cells.Item(43, 11).Value = 622.3333
$ws.Cells.Item(43, 12).Value = 855.1667
$ws.Cells.Item(43, 13).Value = -553.3333
$ws.Cells.Item(43, 14).Value = -993.1667

$ws.Cells.Item(62, 8).Value = 2887.2222
$ws.Cells.Item(62, 9).Value = 2037
$ws.Cells.Item(62, 10).Value = 3950
$ws.Cells.Item(62, 11).Value = 2037
$ws.Cells.Item(62, 12).Value = 3950
$ws.Cells.Item(62, 13).Value = -1413
$ws.Cells.Item(62, 14).Value = -5198

$ws.Cells.Item(65, 8).Value = 2887.2222
$ws.Cells.Item(65, 9).Value = 2037
$ws.Cells.Item(65, 10).Value = 3950
$ws.Cells.Item(65, 11).Value = 10185
$ws.Cells.Item(65, 12).Value = 19750
$ws.Cells.Item(65, 13).Value = -7065
$ws.Cells.Item(65, 14).Value = -25990

$ws.Cells.Item(80, 8).Value = 545
$ws.Cells.Item(80, 9).Value = 200
$ws.Cells.Item(80, 10).Value = 775
$ws.Cells.Item(80, 11).Value = 600
$ws.Cells.Item(80, 12).Value = 2325
$ws.Cells.Item(80, 13).Value = 398
$ws.Cells.Item(80, 14).Value = -4321

$ws.Cells.Item(83, 8).Value = 545
$ws.Cells.Item(83, 9).Value = 200
$ws.Cells.Item(83, 10).Value = 775
$ws.Cells.Item(83, 11).Value = 1800
$ws.Cells.Item(83, 12).Value = 6975
$ws.Cells.Item(83, 13).Value = 3192
$ws.Cells.Item(83, 14).Value = -16959

$ws.Cells.Item(86, 8).Value = 100984.82
$ws.Cells.Item(86, 9).Value = 100984.82
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 100984.82
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = -99861.82000000001
$ws.Cells.Item(86, 14).Value = ""

$ws.Cells.Item(89, 8).Value = 100984.82
$ws.Cells.Item(89, 9).Value = 100984.82
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 504924.1
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = -499308.1
$ws.Cells.Item(89, 14).Value = ""

$ws.Cells.Item(125, 8).Value = 901.5
$ws.Cells.Item(125, 9).Value = 800
$ws.Cells.Item(125, 10).Value = 1003
$ws.Cells.Item(125, 11).Value = 7200
$ws.Cells.Item(125, 12).Value = 9027
$ws.Cells.Item(125, 13).Value = -4740
$ws.Cells.Item(125, 14).Value = -13947

$ws.Cells.Item(127, 8).Value = 1024.8334
$ws.Cells.Item(127, 9).Value = 239.5
$ws.Cells.Item(127, 11).Value = 718.5
$ws.Cells.Item(127, 13).Value = 4241.5

$ws.Cells.Item(129, 8).Value = 1936.3704
$ws.Cells.Item(129, 9).Value = 667.7778
$ws.Cells.Item(129, 10).Value = 2570.6667
$ws.Cells.Item(129, 11).Value = 2003.3334
$ws.Cells.Item(129, 12).Value = 7712.000100000001
$ws.Cells.Item(129, 13).Value = 2996.6666
$ws.Cells.Item(129, 14).Value = -17712.0001

$ws.Cells.Item(131, 8).Value = 4304.7144
$ws.Cells.Item(131, 9).Value = 790
$ws.Cells.Item(131, 10).Value = 5131.706
$ws.Cells.Item(131, 11).Value = 2370
$ws.Cells.Item(131, 12).Value = 15395.118
$ws.Cells.Item(131, 13).Value = 2670
$ws.Cells.Item(131, 14).Value = -25475.118

$ws.Cells.Item(137, 8).Value = 1595.6666
$ws.Cells.Item(137, 9).Value = 2162.8572
$ws.Cells.Item(137, 10).Value = 1234.7273
$ws.Cells.Item(137, 11).Value = 6488.571599999999
$ws.Cells.Item(137, 12).Value = 3704.1819
$ws.Cells.Item(137, 13).Value = -3938.571599999999
$ws.Cells.Item(137, 14).Value = -8804.1819

$ws.Cells.Item(138, 8).Value = 2765.5823
$ws.Cells.Item(138, 9).Value = 2189.1428
$ws.Cells.Item(138, 10).Value = 2974.2932
$ws.Cells.Item(138, 11).Value = 6567.428400000001
$ws.Cells.Item(138, 12).Value = 8922.8796
$ws.Cells.Item(138, 13).Value = -1427.428400000001
$ws.Cells.Item(138, 14).Value = -19202.8796

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(16, 8).Value = 2717.2
$ws.Cells.Item(16, 9).Value = 2146.5
$ws.Cells.Item(16, 11).Value = 2146.5
$ws.Cells.Item(16, 13).Value = -1859.5

$ws.Cells.Item(32, 8).Value = 7370.72
$ws.Cells.Item(32, 9).Value = 6103.5493
$ws.Cells.Item(32, 10).Value = 20183.223
$ws.Cells.Item(32, 11).Value = 6103.5493
$ws.Cells.Item(32, 12).Value = 20183.223
$ws.Cells.Item(32, 13).Value = -5816.5493
$ws.Cells.Item(32, 14).Value = -20757.223

$ws.Cells.Item(57, 8).Value = 26000
$ws.Cells.Item(57, 9).Value = 26000
$ws.Cells.Item(57, 11).Value = 26000
$ws.Cells.Item(57, 13).Value = -25516

$ws.Cells.Item(74, 8).Value = 1324.1613
$ws.Cells.Item(74, 9).Value = 1163.0435
$ws.Cells.Item(74, 11).Value = 1163.0435
$ws.Cells.Item(74, 13).Value = -289.0435

$ws.Cells.Item(77, 8).Value = 1324.1613
$ws.Cells.Item(77, 9).Value = 1163.0435
$ws.Cells.Item(77, 11).Value = 5815.2175
$ws.Cells.Item(77, 13).Value = -1447.2175

$ws.Cells.Item(101, 8).Value = 17734.666
$ws.Cells.Item(101, 10).Value = 17734.666
$ws.Cells.Item(101, 12).Value = 17734.666
$ws.Cells.Item(101, 14).Value = -24224.666

$ws.Cells.Item(122, 8).Value = 1394.3334
$ws.Cells.Item(122, 9).Value = 1380.2667
$ws.Cells.Item(122, 10).Value = 1464.6666
$ws.Cells.Item(122, 11).Value = 4140.800099999999
$ws.Cells.Item(122, 12).Value = 4393.9998
$ws.Cells.Item(122, 13).Value = -1690.800099999999
$ws.Cells.Item(122, 14).Value = -9293.9998

$ws.Cells.Item(126, 8).Value = 5333.3335
$ws.Cells.Item(126, 9).Value = 5333.3335
$ws.Cells.Item(126, 11).Value = 16000.0005
$ws.Cells.Item(126, 13).Value = -13530.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 579
$ws.Cells.Item(11, 10).Value = 605
$ws.Cells.Item(11, 12).Value = 605
$ws.Cells.Item(11, 14).Value = -885

$ws.Cells.Item(109, 8).Value = 31130
$ws.Cells.Item(109, 10).Value = 31130
$ws.Cells.Item(109, 12).Value = 31130
$ws.Cells.Item(109, 14).Value = -33904

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(43, 8).Value = 20995.223
$ws.Cells.Item(43, 10).Value = 20995.223
$ws.Cells.Item(43, 12).Value = 20995.223
$ws.Cells.Item(43, 14).Value = -21363.223

$ws.Cells.Item(101, 8).Value = 20995.223
$ws.Cells.Item(101, 10).Value = 20995.223
$ws.Cells.Item(101, 12).Value = 20995.223
$ws.Cells.Item(101, 14).Value = -27485.223

$ws.Cells.Item(131, 8).Value = 21526
$ws.Cells.Item(131, 10).Value = 21526
$ws.Cells.Item(131, 12).Value = 21526
$ws.Cells.Item(131, 14).Value = -31606

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 662.7
$ws.Cells.Item(5, 9).Value = 400.08
$ws.Cells.Item(5, 11).Value = 1200.24
$ws.Cells.Item(5, 13).Value = -1088.24

$ws.Cells.Item(20, 8).Value = 4715.385
$ws.Cells.Item(20, 10).Value = 4715.385
$ws.Cells.Item(20, 12).Value = 14146.155
$ws.Cells.Item(20, 14).Value = -14600.155

$ws.Cells.Item(135, 8).Value = 662.7
$ws.Cells.Item(135, 9).Value = 400.08
$ws.Cells.Item(135, 11).Value = 3600.72
$ws.Cells.Item(135, 13).Value = -1065.72

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 5887.7144
$ws.Cells.Item(122, 9).Value = 6242.8
$ws.Cells.Item(122, 11).Value = 18728.4
$ws.Cells.Item(122, 13).Value = -16278.4

$ws.Cells.Item(123, 8).Value = 19732.5
$ws.Cells.Item(123, 10).Value = 19732.5
$ws.Cells.Item(123, 12).Value = 19732.5
$ws.Cells.Item(123, 14).Value = -24632.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(103, 8).Value = 46320.4
$ws.Cells.Item(103, 10).Value = 46320.4
$ws.Cells.Item(103, 12).Value = 46320.4
$ws.Cells.Item(103, 14).Value = -48664.4

$ws.Cells.Item(122, 8).Value = 16273.125
$ws.Cells.Item(122, 9).Value = 22837.4
$ws.Cells.Item(122, 10).Value = 5332.6665
$ws.Cells.Item(122, 11).Value = 68512.20000000001
$ws.Cells.Item(122, 12).Value = 15997.9995
$ws.Cells.Item(122, 13).Value = -66062.20000000001
$ws.Cells.Item(122, 14).Value = -20897.9995

$ws.Cells.Item(136, 8).Value = 1799.0605
$ws.Cells.Item(136, 9).Value = 1531.2593
$ws.Cells.Item(136, 10).Value = 3004.1667
$ws.Cells.Item(136, 11).Value = 4593.7779
$ws.Cells.Item(136, 12).Value = 9012.500100000001
$ws.Cells.Item(136, 13).Value = -2043.7779
$ws.Cells.Item(136, 14).Value = -14112.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(92, 8).Value = 30550
$ws.Cells.Item(92, 10).Value = 30550
$ws.Cells.Item(92, 12).Value = 30550
$ws.Cells.Item(92, 14).Value = -35542

$ws.Cells.Item(132, 8).Value = 1962.4889
$ws.Cells.Item(132, 9).Value = 1614.7391
$ws.Cells.Item(132, 10).Value = 2326.0454
$ws.Cells.Item(132, 11).Value = 4844.2173
$ws.Cells.Item(132, 12).Value = 6978.1362
$ws.Cells.Item(132, 13).Value = -2314.2173
$ws.Cells.Item(132, 14).Value = -12038.1362

$ws.Cells.Item(136, 8).Value = 2362.0857
$ws.Cells.Item(136, 9).Value = 2344.2903
$ws.Cells.Item(136, 10).Value = 2500
$ws.Cells.Item(136, 11).Value = 7032.8709
$ws.Cells.Item(136, 12).Value = 7500
$ws.Cells.Item(136, 13).Value = -4482.8709
$ws.Cells.Item(136, 14).Value = -12600

